# Update by John on Task 5.1
#
# - Sheet1: row 14 (Task 5.1) gets its Test Result / Timestamp filled in.
# - Summary: roll-up numbers refreshed for the new Pass result.
# - "Task ID 5": a new "Task 5.1" subtask block (with its own screenshot)
#   is appended below the existing "Task 5.0" block, and the Task 5.0
#   screenshot gets a solid outline to match the other tabs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 - master task table: Task 5.1 (row 14) now has a result
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Reuse the formatting already used for the Task 5.0 "Pass" cell (F13)
# so the green fill / wrap style carries over to F14.
$ws1.Range("F13").Copy()
$ws1.Range("F14").PasteSpecial(-4122)

$ws1.Range("F14").Value = "Pass"
$ws1.Range("G14").Value = "2025-05-04 16:52:19"

# ---------------------------------------------------------------------
# Summary - recompute the headline stats now that Task 5.1 has a result
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 4
$wsSummary.Range("B5").Value = "30.77%"
$wsSummary.Range("B6").Value = "5.1"
$wsSummary.Range("B8").Value = "2025-05-04 16:52:19"
$wsSummary.Range("B12").Value = "[█████████████-------] 69%"
$wsSummary.Range("B20").Value = 9

# John now has 2 tests logged, bumping him above Yash/Paul in the
# "tests per tester" leaderboard.
$wsSummary.Range("A42").Value = "John"
$wsSummary.Range("B42").Value = 2
$wsSummary.Range("A43").Value = "Yash"
$wsSummary.Range("B43").Value = 1
$wsSummary.Range("A44").Value = "Paul"
$wsSummary.Range("B44").Value = 1

# ---------------------------------------------------------------------
# Task ID 5 - append the "Task 5.1" subtask block
# ---------------------------------------------------------------------
$wsTask = $wb.Worksheets.Item("Task ID 5")

# Clone the layout/formatting of the Task 5.0 info block (rows 1-4) and
# its Test Result block (rows 20-21) into the new block's rows.
$wsTask.Range("A1:B4").Copy()
$wsTask.Range("A23:B26").PasteSpecial(-4122)

$wsTask.Range("A20:B21").Copy()
$wsTask.Range("A42:B43").PasteSpecial(-4122)

$wsTask.Range("A23").Value = "Subtask"
$wsTask.Range("B23").Value = "Task 5.1"
$wsTask.Range("A24").Value = "Navigation"
$wsTask.Range("B24").Value = "Settings > Notifications > Opt-In"
$wsTask.Range("A25").Value = "Tester Name"
$wsTask.Range("B25").Value = "John"
$wsTask.Range("A26").Value = "Timestamp"
$wsTask.Range("B26").Value = "2025-05-04 16:52:19"

# Row 27 is the spacer row that holds the new screenshot, same as row 5
# does for the Task 5.0 screenshot.
$wsTask.Rows.Item(27).RowHeight = 100

$wsTask.Range("A42").Value = "Test Result"
$wsTask.Range("B42").Value = "Pass"
$wsTask.Range("A43").Value = "Comment"
$wsTask.Range("B43").Value = "test"

# Give the existing Task 5.0 screenshot a solid outline, matching the
# other "Task ID" tabs that already carry this styling.
$pic1 = $wsTask.Shapes.Item(1)
$pic1.Line.Visible = $true
$pic1.Line.DashStyle = 1

# Attach the new Task 5.1 screenshot under the new subtask block, sized
# and positioned the same way as the Task 5.0 screenshot.
$anchorCell = $wsTask.Range("A27")
$wsTask.Shapes.AddPicture("/tmp/work/local_image.png", 0, -1, $anchorCell.Left, $anchorCell.Top, 450, 300) | Out-Null
